$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("scale")
$ws.Activate()
